# DDAf_2023_Tableau_annexe_Tab19.xlsx - "Add files via upload"
# The authors re-exported this annex table from their data pipeline; almost every
# change in the underlying OOXML is a recalculated statistic in column R
# ("Valeur totale des exportations de biens et services en % du PIB, 2021") plus a
# handful of other recomputed cells (H67, L67, N70, and the full rows 97-98 ranges),
# and one cell (R26, Érythrée) whose figure became unavailable ("..") in this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab19")

# --- Column R ("Valeur totale des exportations de biens et services en % du PIB, 2021") ---
# bulk of the refreshed figures, row number -> new value
$colR = @{
    3 = 44.983369354109797
    4 = 41.960570590450502
    5 = 45.037318436425998
    6 = 44.4455922444536
    7 = 10.3622362790511
    8 = 39.8111704972711
    9 = 33.957445076742701
    10 = 31.111990655453798
    11 = 52.950702242066598
    12 = 18.268420102739999
    13 = 34.502451308764499
    14 = 8.12233191914566
    15 = 16.441806010866902
    16 = 13.692146062069099
    17 = 34.559806862069998
    18 = 66.052201281445903
    19 = 38.993071731061399
    20 = 44.624162319167297
    21 = 39.888258732818699
    22 = 15.3632638728612
    23 = 32.666933938057603
    24 = 11.5195131831656
    25 = 149.33312334600899
    27 = 8.5738621328594995
    28 = 10.7623285479165
    29 = 23.0317273620702
    30 = 27.838981225178699
    31 = 19.088529634091898
    32 = 116.66493004945001
    33 = 13.3178343668908
    34 = 71.501738670582995
    35 = 16.993018195085298
    36 = 14.175209086724699
    37 = 14.4746086054706
    38 = 14.441826690928799
    39 = 25.605436625959499
    40 = 10.5531774155445
    42 = 33.205487902222501
    43 = 33.229757648561801
    44 = 42.4403651619307
    45 = 18.816926072228799
    46 = 23.4835027649777
    47 = 29.457614399860699
    48 = 29.639481175743299
    49 = 22.589718820485
    50 = 9.1314405670974601
    51 = 30.027140660352998
    52 = 64.805956810950704
    53 = 19.428377144420999
    54 = 32.0112993234793
    55 = 27.222714232390601
    56 = 15.4966004375969
    57 = 11.4996686864061
    58 = 23.420772400342599
    59 = 26.863274727982201
    60 = 22.8434868611367
    61 = 17.674116448833999
    62 = 21.380587503797599
    63 = 29.0241195092248
    64 = 26.8460489152917
    65 = 24.090908250574699
    66 = 28.639168740902601
    67 = 15.7247606941299
    68 = 17.539831637828499
    69 = 17.4506879069239
    70 = 36.434648860728501
    71 = 17.674116448833999
    72 = 13.1362761632726
    73 = 32.338526904563103
    74 = 30.575382782963299
    75 = 43.3975442533836
    76 = 52.985083071879501
    77 = 21.219838907156699
    78 = 50.201959459432601
    79 = 29.930302691805601
    80 = 33.437073793099799
    81 = 37.650932556709002
    82 = 19.6853192705649
    83 = 28.2433343612032
    84 = 22.672854361169701
    85 = 10.6337358712987
    86 = 18.764616357292901
    87 = 25.209937138664401
    88 = 32.185209258026298
    89 = 26.197198144810798
    90 = 32.422092984809197
    91 = 25.195366261229101
    92 = 15.3101750673393
    93 = 32.281796600955701
    94 = 119.611789619783
    95 = 20.3560475644985
    96 = 31.229914599945701
}
foreach ($row in $colR.Keys) {
    $ws.Range("R$row").Value = $colR[$row]
}

# Érythrée (row 26): figure no longer available in this refresh
$ws.Range("R26").Value = ".."

# A few other individually recomputed cells
$ws.Range("H67").Value = 5.2917556161347301
$ws.Range("L67").Value = 18.601144175363899
$ws.Range("N70").Value = 3.00757906215029

# Row 97 (full row of refreshed figures)
$row97 = @{
    "C" = 9.2329908505254696
    "D" = 7.11636757045701
    "E" = 11.370953785518401
    "F" = 15.8363945094399
    "G" = 24.4636653418919
    "H" = 4.4669461356271798
    "I" = 11.807813104587201
    "J" = 9.7750280506337699
    "K" = 13.507604281606
    "L" = 25.983298779366599
    "M" = 18.089948175338598
    "N" = 3.6910943397805802
    "O" = 21583.432691000002
    "P" = 96378.972787999999
    "Q" = 18.296873994183098
    "R" = 18.551904407748001
}
foreach ($col in $row97.Keys) {
    $ws.Range("$col" + "97").Value = $row97[$col]
}

# Row 98 (full row of refreshed figures; D98/J98 stay unavailable "..")
$row98 = @{
    "C" = 6.5883938736442298
    "E" = 42.999359961668297
    "F" = 14.546495702032599
    "G" = 23.9815661593687
    "H" = 18.139540819206601
    "I" = 11.9716161625933
    "K" = 74.102027717615599
    "L" = 28.0554149027077
    "M" = 6.4795842201088698
    "N" = 7.4546122748981301
    "O" = 153212.318306
    "P" = 54680.485189000101
    "Q" = 73.697749864480002
    "R" = 19.257656382611401
}
foreach ($col in $row98.Keys) {
    $ws.Range("$col" + "98").Value = $row98[$col]
}
